$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.191.95"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "2.421.08"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.40"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.29"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.70%  "
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.70"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.01%  "
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.87"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "2.852.06"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "60.080.38"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("D17").Value = "2.437.06"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.49"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "327.14"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.74"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.25"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.78"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +3.06%  "
$ws.Range("E28").Value = "  -1.65%  "
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.39"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("E31").Value = "  -2.96%  "
$ws.Range("E32").Value = "  +2.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.404"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.73%  "
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "328.95"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "144.76"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.43%  "
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.09"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0965"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.575"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.03"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("E49").Value = "  -2.38%  "
$ws.Range("E51").Value = "  -0.65%  "
